# "fix name elc ava"
# - Rename the "Availability" sheet to "ELC_AVA".
# - On that sheet, the DKE/DKW availability flags (row 5-7, cols D/E) were
#   missing; the AllRegions flag in C5 was wrong and is cleared, and D/E
#   (DKE/DKW) get the "1" flag for rows 5-7 instead.
# - Leave the final selection on E7, matching the last cell touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Availability")
$ws.Name = "ELC_AVA"

[void]$ws.Range("C5").ClearContents()

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

[void]$ws.Range("E7").Select()
